$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New plain rows 32-42 (default style, same as header A1) ---
$plainValues = @(
    "용인",
    "화성",
    "동두천",
    "평택",
    "의왕",
    "공주",
    "경기",
    "군포",
    "이천",
    "안산",
    "서천군"
)

$row = 32
foreach ($val in $plainValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

# --- New emphasized rows 43-46 (bigger font, left-aligned, taller row) ---
# row starts at 43 here
$emphValues = @("충북", "전남", "전북", "충남")
foreach ($val in $emphValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

# Build the emphasis formatting on the first emphasized cell (A43) ...
$first = $ws.Cells.Item(43, 1)
$first.Font.Size = 12
$first.HorizontalAlignment = -4131
$first.VerticalAlignment = -4108
$ws.Rows.Item(43).RowHeight = 17.25

# ... then copy that formatting onto the remaining emphasized cells so every
# row ends up referencing the very same cell style instead of Excel minting
# a fresh (duplicate) style record per cell.
$first.Copy()
$ws.Range("A44:A46").PasteSpecial(-4122)

$ws.Rows.Item(44).RowHeight = 17.25
$ws.Rows.Item(45).RowHeight = 17.25
$ws.Rows.Item(46).RowHeight = 17.25

# Match the author's final on-screen selection.
[void]$ws.Range("E16").Select()

"done"
